# Updates cryptos list prices / 1h volume changes (GitHub Actions data refresh).
# Also swaps the USDe / RenderToken rows (30 & 31) to reflect their new rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.486.43"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "3.144.10"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.142.08"
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.447"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.61%  "
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.394"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("D13").Value = "3.683.83"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("E14").Value = "  +3.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.29%  "
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "58.484.71"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Value = "3.145.60"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "344.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.95%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.17%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("E39").Value = "  -4.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0674"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.710"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.38%  "
$ws.Range("D44").Value = "3.184.06"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0265"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.56%  "
$ws.Range("D48").Value = "2.298.94"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.69%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.78%  "
